$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push up the boolean value: F2 ("Is Active") should hold an actual
# boolean FALSE instead of the text string "False".
$ws.Range("F2").Value = $false
